# edit.ps1
# Applies the changes described by the target diff:
#  1. Split "jdbc:derby://localhost:1527/" into two runs: "jdbc:derby" and "://localhost:1527/"
#  2. Split ". Username and password: pdc" into ". Username and password: " and "pdc"
#  3. Split " controller" into " " and "controller", making "controller" bold (in addition to
#     the existing italic formatting)
#  4. Trim the long trailing whitespace run down to 17 spaces, then start a new paragraph that
#     reads "The test classes are in the test package." followed by the remaining whitespace
#     (4118 spaces) that used to live at the end of the previous paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force Word to split a contiguous run of text into its own run by
# toggling a character formatting property on and back off again. This keeps
# the visual formatting unchanged but guarantees the sub-range becomes a
# distinct <w:r> element instead of being silently merged back into its
# neighbours.
# ---------------------------------------------------------------------------
function Split-Run($range) {
    $range.Font.Bold = 1
    $range.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 1) jdbc:derby://localhost:1527/  ->  "jdbc:derby" + "://localhost:1527/"
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("jdbc:derby")
$rJdbc = $d.Range($idx, $idx + "jdbc:derby".Length)
Split-Run $rJdbc

# ---------------------------------------------------------------------------
# 2) ". Username and password: pdc"  ->  ". Username and password: " + "pdc"
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("password: pdc") + "password: ".Length
$rPdc = $d.Range($idx, $idx + "pdc".Length)
Split-Run $rPdc

# ---------------------------------------------------------------------------
# 3) " controller"  ->  " " + "controller" (controller becomes bold as well)
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf(" controller")
$rController = $d.Range($idx + 1, $idx + 1 + "controller".Length)
$rController.Font.Bold = 1

# ---------------------------------------------------------------------------
# 4) Trim the trailing run of 4135 spaces to 17 spaces and move the remaining
#    4118 spaces into a new paragraph that begins with
#    "The test classes are in the test package."
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("exit button. ")
$wsStart = $idx + "exit button. ".Length
$splitPos = $wsStart + 17

# Break the paragraph right after the 17th space.
$rBreak = $d.Range($splitPos, $splitPos)
$rBreak.InsertParagraphAfter()

# Insert the new sentence at the very start of the newly created paragraph
# (i.e. right after the paragraph mark we just inserted).
$insertPos = $splitPos + 1
$sentence = "The test classes are in the test package."
$rInsert = $d.Range($insertPos, $insertPos)
$rInsert.InsertBefore($sentence)

# Force the new sentence and the remaining whitespace to live in separate runs.
$rSentence = $d.Range($insertPos, $insertPos + $sentence.Length)
Split-Run $rSentence
